$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the F/G (LinkedIn_Poster / LinkedIn_Posted) values for rows 4 and 5
$ws.Range("F4:G4").ClearContents()
$ws.Range("F5:G5").ClearContents()

# Update the active selection to match the saved state (K7)
$ws.Range("K7").Select()
